$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Feuil2"
$ws3.Range("B10").Value = 2853
$ws3.Range("B11").Value = 2906
$ws3.Range("B12").Value = 2820
$ws3.Range("B13").Value = 2941
$ws3.Range("B14").Value = 2862
$ws3.Range("B15").Formula = "=AVERAGE(B10:B14)"
